# Generate Report for Handback
# Updates the localization-status workbook after a handback run:
#  - Overview sheet: status flips from "Ready for handoff" to
#    "Handed back: in sync with en-US" for both languages.
#  - Per-language sheets (zh-cn, de-de): the "Latest Handback DateTime"
#    is refreshed and the (now stale) "Error Detail" message is cleared
#    because the handback file is in sync again.
#  - The "Status" / "Error Detail" columns are widened / narrowed to fit
#    the new, shorter content.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("K2").Value = "2016-08-30 19:02:19"
$zhcn.Range("P2").Value = ""
$zhcn.Columns.Item(3).ColumnWidth = 29.1666666666667
$zhcn.Columns.Item(16).ColumnWidth = 12.8333333333333

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("K2").Value = "2016-08-30 19:02:27"
$dede.Range("P2").Value = ""
$dede.Columns.Item(3).ColumnWidth = 29.1666666666667
$dede.Columns.Item(16).ColumnWidth = 12.8333333333333

Write-Output "Generate Report for Handback: applied."
